$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# The sheet used to hold a single giant inline-string cell (A1) dumped
# straight from a PDF-to-Excel conversion, styled bold + centered, with
# column A force-widened to fit the whole block of text. We're replacing
# that with a proper grid: one label/value per cell, laid out the way the
# original report is structured (CCM / CPF-CNPJ / month / status / closed
# date header, company name, address, complement/neighborhood/city/state).
#
# Swapping in a brand-new worksheet (same name, same position) gets rid of
# the old custom column width and the old cell's style in one shot, rather
# than trying to claw them back field-by-field.
# ---------------------------------------------------------------------
$oldName = $wb.ActiveSheet.Name
$newWs = $wb.Worksheets.Add()
$wb.Worksheets.Item(2).Delete()
$newWs.Name = $oldName
$ws = $wb.Worksheets.Item(1)

# Restore the page margins the workbook originally had (a fresh sheet
# defaults to different ones).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Labelled cells.
$labels = @{
    "A10" = "CCM:"
    "E10" = "CPF / CNPJ:"
    "I10" = "Mês Referência:"
    "L10" = "Situação:"
    "O10" = "Encerramento:"

    "A12" = "Razão Social:"
    "A13" = "TRIO ARENA UNIVERSITARIA LTDA."

    "A14" = "Endereço:"
    "J14" = "Número:"

    "A16" = "Complemento:"
    "E16" = "Bairro:"
    "I16" = "Cidade:"
    "N16" = "Estado:"
}
foreach ($addr in $labels.Keys) {
    $ws.Range($addr).Value = $labels[$addr]
}

# Value cells directly under each label - left blank (the report engine
# fills these in per-record); they still need to exist as formatted
# cells.
$valueCells = @("A11","E11","I11","L11","O11","A15","J15","A17","E17","I17","N17")
foreach ($addr in $valueCells) {
    $ws.Range($addr).HorizontalAlignment = -4108
}

# Every populated/reserved cell shares the same style: centered both
# ways, regular (non-bold) weight.
$allCells = @($labels.Keys) + $valueCells
foreach ($addr in $allCells) {
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4108
}
